$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H) 14:33:33 -> 14:33:57
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"

$wsZhCn.Range("H4").Value = "2016-08-22 14:33:57"
$wsZhCn.Range("H5").Value = "2016-08-22 14:33:57"
$wsZhCn.Range("H6").Value = "2016-08-22 14:33:57"
$wsZhCn.Range("H7").Value = "2016-08-22 14:33:57"

# de-de sheet: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H) 14:33:38 -> 14:34:11
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"

$wsDeDe.Range("H4").Value = "2016-08-22 14:34:11"
$wsDeDe.Range("H5").Value = "2016-08-22 14:34:11"
$wsDeDe.Range("H6").Value = "2016-08-22 14:34:11"
$wsDeDe.Range("H7").Value = "2016-08-22 14:34:11"

# Overview sheet: rows 4-7 -> Latest HO Xliff Generate Date (G) 14:33:38 -> 14:34:11
# (shares the same text/shared-string as de-de!H4:H7)
$wsOverview.Range("G4").Value = "2016-08-22 14:34:11"
$wsOverview.Range("G5").Value = "2016-08-22 14:34:11"
$wsOverview.Range("G6").Value = "2016-08-22 14:34:11"
$wsOverview.Range("G7").Value = "2016-08-22 14:34:11"
